$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so strings like "1.003" or
# "0.5126" are not silently reinterpreted as numbers (matches the source data,
# which stores these as literal text/inline strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.989.08"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.859.47"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "311.67"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.5126"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "0.08321"
$ws.Range("E9").Value = "  -9.33%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "41.28"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "6.192"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "1.861.70"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "7.176"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "0.00001092"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "90.31"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "0.06624"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "28.017.15"
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "2.250"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").Value = "2.574"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").Value = "2.078.25"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "157.10"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "20.46"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "125.30"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "0.1060"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "1.037"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").Value = "5.584"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "9.559"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "0.06509"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "0.02406"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "0.2150"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "0.6394"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "1.225"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").Value = "11.26"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").Value = "4.853"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "0.6060"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "12.98"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").Value = "1.282"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "3.655"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "1.974"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "1.206"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "120.47"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "79.49"
$ws.Range("E51").Value = "  +1.11%  "

# Restore the default (unstyled) cell style now that the text values are locked in,
# so the workbook formatting matches the original (no explicit style on these cells).
$ws.Range("D2:D51").Style = "Normal"
